$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-25 13:43:59"
$wsZhCn.Range("G3").Value = "2016-01-25 13:44:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-25 13:44:10"
$wsDeDe.Range("G3").Value = "2016-01-25 13:44:59"
